$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 3 (CityIdentifier)
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = 8
$ws.Range("E3").Value = "8 groups of 1 mil each"

# ---------------------------------------------------------------------------
# Row 9 (Census_FirmwareManufacturerIdentifier)
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = 4
$ws.Range("E9").Value = "142, 628, 554, everything else"
$ws.Range("F9").Value = "2.6, 1.2, 1.1, 3.8 mil observations in each group (ordered)"

# ---------------------------------------------------------------------------
# Row 10 (OsBuildLab)
# ---------------------------------------------------------------------------
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = "16299, 17134, everything else"
$ws.Range("F10").Value = "2.6, 3.9, 2.5 mil each (ordered)"

# ---------------------------------------------------------------------------
# Row 11 (Census_OSVersion)
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = "10.0.17134.228, 10.0.17134.165, everything else divided into 6 more groups"
$ws.Range("F11").Value = "Approx 1 mil each"

# ---------------------------------------------------------------------------
# Row 12 (IeVerIdentifier) - also gets the new green highlight fill
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "17134, 16299, everything else"
$ws.Range("F12").Value = "3.9, 2.5, 2.6 mil each (ordered)"

# ---------------------------------------------------------------------------
# Row 13 (GeoNameIdentifier)
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = 3
$ws.Range("E13").Value = "277, then divide the remaining into two equal groups"

# ---------------------------------------------------------------------------
# Row 14 (Census_OSBuildRevision)
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = 8
$ws.Range("E14").Value = "228, then divide the remaining into 7 more equal groups "
$ws.Range("F14").Value = "Each group should have ~ 1.5 mil observations"

# ---------------------------------------------------------------------------
# Row 15 (LocaleEnglishNameIdentifier)
# ---------------------------------------------------------------------------
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = "228, then divide the remaining into 7 more equal groups "
$ws.Range("F15").Value = "Each group should have ~ 1.5 mil observations"

# ---------------------------------------------------------------------------
# Row 16 (CountryIdentifier) - only D and F, no E
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 8
$ws.Range("F16").Value = "Divide into groups of ~ 1 mil"

# ---------------------------------------------------------------------------
# Row 17 (Census_OSBuildNumber)
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = "17134, 16299, everything else"
$ws.Range("F17").Value = "4, 2.4, ~3 mil each"

# ---------------------------------------------------------------------------
# Row 18 (Census_OSUILocaleIdentifier)
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = "31, 34, divide the rest into 3 groups"
$ws.Range("F18").Value = "3.1, 1, ~1+ mil each"

# ---------------------------------------------------------------------------
# Row 19 (AppVersion) - whole row gets the border style (style 1), incl a new
# empty bordered F19 cell
# ---------------------------------------------------------------------------
$ws.Range("A19:F19").Borders(9).LineStyle = 1
$ws.Range("A19:F19").Borders(9).Weight = 2

# ---------------------------------------------------------------------------
# Row 20 (OsBuild) - also gets the new green highlight fill; E20 keeps the
# header-like style (style 2)
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = "17134, 16299, everything else"
$ws.Range("F20").Value = "3.9, 2.5, 2.6 mil each (ordered)"

# ---------------------------------------------------------------------------
# Row 21 (EngineVersion)
# ---------------------------------------------------------------------------
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = "1.1.15200.1, 1.1.15100.1, everything else"
$ws.Range("F21").Value = "3.8, 3.6, ~1 mil each"

# ---------------------------------------------------------------------------
# Row 22 (OsVer)
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = "10.0.0.0, everything else"
$ws.Range("F22").Value = "8.6, <1"

# ---------------------------------------------------------------------------
# Row 23 (Census_ChassisTypeName)
# ---------------------------------------------------------------------------
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = "notebook, everything else"
$ws.Range("F23").Value = "5.2 mil are notebook. Everything else could be grouped into another massive bundle or it could be sorted better. "

# ---------------------------------------------------------------------------
# Row 24 (OrganizationIdentifier) - whole row gets the border style (style 1)
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = "27, 18,  everything else"
$ws.Range("F24").Value = "4.1, 1.7, ~3 mil each"
$ws.Range("A24:F24").Borders(9).LineStyle = 1
$ws.Range("A24:F24").Borders(9).Weight = 2

# ---------------------------------------------------------------------------
# Rows 25-48: copy the unique-value count into the "new value count" column,
# i.e. each of these low-cardinality features is kept as-is (no grouping)
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 33
$ws.Range("D26").Value = 32
$ws.Range("D27").Value = 30
$ws.Range("D28").Value = 22
$ws.Range("D29").Value = 16
$ws.Range("D30").Value = 14
$ws.Range("D31").Value = 13
$ws.Range("D32").Value = 11
$ws.Range("D33").Value = 10
$ws.Range("D34").Value = 9
$ws.Range("D35").Value = 9
$ws.Range("D36").Value = 9
$ws.Range("D37").Value = 8
$ws.Range("D38").Value = 8
$ws.Range("D39").Value = 8
$ws.Range("D40").Value = 7
$ws.Range("D41").Value = 6
$ws.Range("D42").Value = 6
$ws.Range("D43").Value = 5
$ws.Range("D44").Value = 5
$ws.Range("D45").Value = 4
$ws.Range("D46").Value = 3
$ws.Range("D47").Value = 3
$ws.Range("D48").Value = 3

# Row 49 (HasDetections) already has the border style; it just gains a value
$ws.Range("D49").Value = 2

# ---------------------------------------------------------------------------
# Highlight fill (theme "Green, Accent 6, Lighter 80%") for the two rows
# (IeVerIdentifier / OsBuild) that are duplicates of another row's grouping
# ---------------------------------------------------------------------------
$ws.Range("B12").Interior.ThemeColor = 10
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# E20 keeps the same "header" style used elsewhere in column F / row 50
$ws.Range("E20").Copy() | Out-Null
$ws.Range("F1").Copy() | Out-Null

# ---------------------------------------------------------------------------
# View state: scroll position & selection
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("G28").Select() | Out-Null

# ---------------------------------------------------------------------------
# Column E width grew to fit the longer group descriptions
# ---------------------------------------------------------------------------
$ws.Columns("E:E").ColumnWidth = 63.53

Write-Host "done"
